# Update yearly financials for AEMD: Minority Interest, Income After Tax,
# Net Income From Continuing Ops, and Effect Of Accounting Changes rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AEMD")

# Row 20 - Minority Interest
$ws.Range("D20").Value = -100
$ws.Range("E20").Value = -300
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 200
$ws.Range("I20").Value = 0

# Row 21 - Income After Tax
$ws.Range("D21").Value = -5300
$ws.Range("E21").Value = -7000
$ws.Range("F21").Value = -4300
$ws.Range("G21").Value = -6500
$ws.Range("I21").Value = -3700

# Row 22 - Net Income From Continuing Ops
$ws.Range("D22").Value = 400
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 600
$ws.Range("G22").Value = 400
$ws.Range("I22").Value = 1100

# Row 32 - Effect Of Accounting Changes
$ws.Range("D32").Value = 100
$ws.Range("E32").Value = 300
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = -200
$ws.Range("I32").Value = 0
